$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing answer for the second question
$ws.Range("B2").Value = "Rose-Hulman is ranked number 1"

# Column B now needs its own (wider) width, separate from column C
$ws.Range("B1").ColumnWidth = 28.6665

# Move the active selection to C5
$ws.Range("C5").Select()
